# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to match the newly generated gh-pages output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of row -> new F-column value for worksheet "展览"
$sheetExhibition = $wb.Worksheets.Item("展览")
$exhibitionUpdates = @{
    2  = 280
    3  = 1160
    4  = 16577
    5  = 18
    6  = 1624
    7  = 58
    8  = 357
    9  = 203
    11 = 11540
    12 = 24
    13 = 1207
    14 = 4566
    15 = 399
    18 = 870
    20 = 146
    21 = 5214
}
foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Map of row -> new F-column value for worksheet "全部类型"
$sheetAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    2  = 280
    4  = 1160
    5  = 16577
    6  = 18
    7  = 1624
    8  = 58
    9  = 357
    10 = 203
    14 = 11540
    15 = 24
    16 = 1207
    17 = 4566
    18 = 399
    21 = 870
    23 = 146
    24 = 5214
}
foreach ($row in $allUpdates.Keys) {
    $sheetAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
